$wb = $excel.ActiveWorkbook

# --- Logs sheet: append a new log row (row 12) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Retour status"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("D12").Value = "Retour / Terugbetaling"
$logs.Range("F12").Value = "2025-08-26 21:31:25"
$logs.Range("G12").Value = "Nee"
$logs.Range("H12").Value = "Ja"
$logs.Range("I12").Value = "Nee"
$logs.Range("J12").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too.
$logs.Range("D2:D11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D12"))
$logs.Range("G2:G11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G12"))
$logs.Range("H2:H11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H12"))
$logs.Range("I2:I11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I12"))
$logs.Range("J2:J11").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J12"))

# --- Dashboard sheet: bump the "Retour / Terugbetaling" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 7
